$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with freshly scraped values.
# NumberFormat is forced to text ("@") before writing the Price column so that
# values such as "67.950.05" or "601.64" are stored as literal strings rather
# than being auto-parsed into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.950.05'
$ws.Range("E2").Value = '  -2.00%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.804.05'
$ws.Range("E3").Value = '  +1.00%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.64'
$ws.Range("E5").Value = '  -2.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.31'
$ws.Range("E6").Value = '  -3.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.804.93'
$ws.Range("E7").Value = '  +1.09%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("E10").Value = '  -4.13%  '
$ws.Range("E11").Value = '  -6.14%  '
$ws.Range("E12").Value = '  -3.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.88'
$ws.Range("E13").Value = '  -2.78%  '
$ws.Range("E14").Value = '  -3.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.438.44'
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.803.27'
$ws.Range("E16").Value = '  +0.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.912.86'
$ws.Range("E17").Value = '  -2.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.28'
$ws.Range("E18").Value = '  -3.60%  '
$ws.Range("E19").Value = '  -3.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.26'
$ws.Range("E20").Value = '  +5.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '494.78'
$ws.Range("E21").Value = '  -2.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.25'
$ws.Range("E22").Value = '  -1.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.86'
$ws.Range("E24").Value = '  -0.82%  '
$ws.Range("E25").Value = '  -4.58%  '
$ws.Range("E26").Value = '  +7.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.38'
$ws.Range("E28").Value = '  -3.78%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("E31").Value = '  -3.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.08'
$ws.Range("E32").Value = '  +7.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.86'
$ws.Range("E33").Value = '  -1.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  -3.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.84'
$ws.Range("E37").Value = '  -4.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.333'
$ws.Range("E38").Value = '  -2.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '462.96'
$ws.Range("E39").Value = '  +1.48%  '
$ws.Range("E40").Value = '  -5.43%  '
$ws.Range("E41").Value = '  -2.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '49.09'
$ws.Range("E42").Value = '  -1.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.85'
$ws.Range("E43").Value = '  -3.93%  '
$ws.Range("E44").Value = '  -1.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.96'
$ws.Range("E45").Value = '  -8.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.846.06'
$ws.Range("E47").Value = '  -3.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '139.75'
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("E49").Value = '  -2.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.04'
$ws.Range("E50").Value = '  +16.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.95'
$ws.Range("E51").Value = '  -4.96%  '
